$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 223-227 with refreshed odds data ---

# Row 223
$ws.Range("B223").Value = 6978388
$ws.Range("E223").Value = 45388.48958333334
$ws.Range("F223").Value = "FC Hebar Pazardzhik"
$ws.Range("G223").Value = "Etar 1924 Veliko Tarnovo"
$ws.Range("K223").Value = 1.75
$ws.Range("L223").Value = 3.5
$ws.Range("M223").Value = 4.75
$ws.Range("N223").Value = 1.8
$ws.Range("O223").Value = 3.4
$ws.Range("P223").Value = 5
$ws.Range("Q223").Value = -0.75
$ws.Range("R223").Value = 2.05
$ws.Range("S223").Value = 1.8
$ws.Range("T223").Value = 2
$ws.Range("U223").Value = 1.825
$ws.Range("V223").Value = 2.025

# Row 224
$ws.Range("B224").Value = 6978438
$ws.Range("E224").Value = 45388.59375
$ws.Range("F224").Value = "Arda Kardzhali"
$ws.Range("G224").Value = "Slavia Sofia"
$ws.Range("K224").Value = 2.05
$ws.Range("L224").Value = 3.2
$ws.Range("M224").Value = 3.75
$ws.Range("N224").Value = 2.25
$ws.Range("O224").Value = 3.1
$ws.Range("P224").Value = 3.5
$ws.Range("Q224").Value = -0.25
$ws.Range("R224").Value = 1.95
$ws.Range("S224").Value = 1.9
$ws.Range("T224").Value = 2
$ws.Range("U224").Value = 1.875
$ws.Range("V224").Value = 1.975

# Row 225
$ws.Range("B225").Value = 8035921
$ws.Range("E225").Value = 45389.36458333334
$ws.Range("F225").Value = "Levski Sofia"
$ws.Range("G225").Value = "CSKA Sofia"
$ws.Range("K225").Value = 2.8
$ws.Range("L225").Value = 3.1
$ws.Range("M225").Value = 2.6
$ws.Range("N225").Value = 2.8
$ws.Range("O225").Value = 3.1
$ws.Range("P225").Value = 2.7
$ws.Range("Q225").Value = 0
$ws.Range("R225").Value = 1.975
$ws.Range("S225").Value = 1.875
$ws.Range("T225").Value = 2
$ws.Range("U225").Value = 2.1
$ws.Range("V225").Value = 1.775

# Row 226
$ws.Range("B226").Value = 8035920
$ws.Range("E226").Value = 45389.48958333334
$ws.Range("F226").Value = "Botev Plovdiv"
$ws.Range("G226").Value = "CSKA 1948 Sofia"
$ws.Range("K226").Value = 2.1
$ws.Range("L226").Value = 3.2
$ws.Range("M226").Value = 3.6
$ws.Range("N226").Value = 1.909
$ws.Range("O226").Value = 3.3
$ws.Range("P226").Value = 4.5
$ws.Range("Q226").Value = -0.5
$ws.Range("R226").Value = 1.925
$ws.Range("S226").Value = 1.925
$ws.Range("T226").Value = 2.25
$ws.Range("U226").Value = 1.825
$ws.Range("V226").Value = 2.025

# Row 227
$ws.Range("B227").Value = 6978446
$ws.Range("E227").Value = 45389.59375
$ws.Range("F227").Value = "Pirin Blagoevgrad"
$ws.Range("G227").Value = "Ludogorets Razgrad"
$ws.Range("K227").Value = 15
$ws.Range("L227").Value = 6
$ws.Range("M227").Value = 1.2
$ws.Range("N227").Value = 19
$ws.Range("O227").Value = 7
$ws.Range("P227").Value = 1.166
$ws.Range("Q227").Value = 2
$ws.Range("R227").Value = 2
$ws.Range("S227").Value = 1.85
$ws.Range("T227").Value = 3
$ws.Range("U227").Value = 1.95
$ws.Range("V227").Value = 1.9

# --- Append new rows 228-230 for newly played/scheduled matches ---

# Row 228
$ws.Range("A227").Copy()
$ws.Range("A228").PasteSpecial(-4122)
$ws.Range("E227").Copy()
$ws.Range("E228").PasteSpecial(-4122)
$ws.Range("A228").Value = 226
$ws.Range("B228").Value = 6978437
$ws.Range("C228").Value = "Bulgaria First League"
$ws.Range("D228").Value = "Bulgaria First League"
$ws.Range("E228").Value = 45390.45833333334
$ws.Range("F228").Value = "Botev Vratsa"
$ws.Range("G228").Value = "Krumovgrad"
$ws.Range("K228").Value = 3
$ws.Range("L228").Value = 3.2
$ws.Range("M228").Value = 2.4
$ws.Range("N228").Value = 3
$ws.Range("O228").Value = 3.1
$ws.Range("P228").Value = 2.55
$ws.Range("Q228").Value = 0
$ws.Range("R228").Value = 2.1
$ws.Range("S228").Value = 1.775
$ws.Range("T228").Value = 2
$ws.Range("U228").Value = 1.9
$ws.Range("V228").Value = 1.95
$ws.Range("W228").Value = 0
$ws.Range("X228").Value = 0
$ws.Range("Y228").Value = 0
$ws.Range("Z228").Value = 0
$ws.Range("AA228").Value = 0

# Row 229
$ws.Range("A228").Copy()
$ws.Range("A229").PasteSpecial(-4122)
$ws.Range("E228").Copy()
$ws.Range("E229").PasteSpecial(-4122)
$ws.Range("A229").Value = 227
$ws.Range("B229").Value = 6978448
$ws.Range("C229").Value = "Bulgaria First League"
$ws.Range("D229").Value = "Bulgaria First League"
$ws.Range("E229").Value = 45390.5625
$ws.Range("F229").Value = "Lokomotiv Plovdiv"
$ws.Range("G229").Value = "Cherno More Varna"
$ws.Range("K229").Value = 2.9
$ws.Range("L229").Value = 3.1
$ws.Range("M229").Value = 2.5
$ws.Range("N229").Value = 3
$ws.Range("O229").Value = 3.1
$ws.Range("P229").Value = 2.5
$ws.Range("Q229").Value = 0
$ws.Range("R229").Value = 2.1
$ws.Range("S229").Value = 1.775
$ws.Range("T229").Value = 2.25
$ws.Range("U229").Value = 2.05
$ws.Range("V229").Value = 1.8
$ws.Range("W229").Value = 0
$ws.Range("X229").Value = 0
$ws.Range("Y229").Value = 0
$ws.Range("Z229").Value = 0
$ws.Range("AA229").Value = 0

# Row 230
$ws.Range("A229").Copy()
$ws.Range("A230").PasteSpecial(-4122)
$ws.Range("E229").Copy()
$ws.Range("E230").PasteSpecial(-4122)
$ws.Range("A230").Value = 228
$ws.Range("B230").Value = 6978449
$ws.Range("C230").Value = "Bulgaria First League"
$ws.Range("D230").Value = "Bulgaria First League"
$ws.Range("E230").Value = 45391.52083333334
$ws.Range("F230").Value = "Lokomotiv 1929 Sofia"
$ws.Range("G230").Value = "Beroe"
$ws.Range("K230").Value = 2.25
$ws.Range("L230").Value = 3
$ws.Range("M230").Value = 3.5
$ws.Range("N230").Value = 2.25
$ws.Range("O230").Value = 3.1
$ws.Range("P230").Value = 3.5
$ws.Range("Q230").Value = -0.25
$ws.Range("R230").Value = 1.925
$ws.Range("S230").Value = 1.925
$ws.Range("T230").Value = 2.25
$ws.Range("U230").Value = 2.025
$ws.Range("V230").Value = 1.825
$ws.Range("W230").Value = 0
$ws.Range("X230").Value = 0
$ws.Range("Y230").Value = 0
$ws.Range("Z230").Value = 0
$ws.Range("AA230").Value = 0
